$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new time-tracking entries appended after the existing A2:B16 data block
# (row 17: 23/10/2013 - 1h20; row 18: 25/10/2013 - 1h55), following the same
# "Data" / "Quantidade de horas" pattern as the rows above them.

# Row 17
$ws.Range("A17").Value = 41570
$ws.Range("B17").Value = 0.055555555555555552
$ws.Range("A16:B16").Copy()
$ws.Range("A17:B17").PasteSpecial(-4122)

# Row 18
$ws.Range("A18").Value = 41572
$ws.Range("B18").Value = 0.079861111111111105
$ws.Range("A16:B16").Copy()
$ws.Range("A18:B18").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Selection moves to the next empty cell in column B, mirroring the cursor
# position left behind after data entry.
$ws.Range("B19").Select()
